$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 590.65216
$ws.Range("I2").Value = 219.46666
$ws.Range("K2").Value = 219.46666
$ws.Range("M2").Value = -106.46666
$ws.Range("H17").Value = 1868.2
$ws.Range("J17").Value = 1868.2
$ws.Range("L17").Value = 5604.6
$ws.Range("N17").Value = -5940.6
$ws.Range("H80").Value = 3597.3333
$ws.Range("J80").Value = 2623.75
$ws.Range("L80").Value = 7871.25
$ws.Range("N80").Value = -9867.25
$ws.Range("H83").Value = 3597.3333
$ws.Range("J83").Value = 2623.75
$ws.Range("L83").Value = 23613.75
$ws.Range("N83").Value = -33597.75
$ws.Range("H88").Value = 1384.25
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 1384.25
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = ""
$ws.Range("M88").Value = 1384.25
$ws.Range("N88").Value = -2196.25
$ws.Range("H91").Value = 1384.25
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1384.25
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = ""
$ws.Range("M91").Value = 1384.25
$ws.Range("N91").Value = -4192.25
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = ""
$ws.Range("M111").Value = ""
$ws.Range("N111").Value = 0
$ws.Range("H138").Value = 4385.3613
$ws.Range("I138").Value = 2587.6
$ws.Range("J138").Value = 5076.8076
$ws.Range("K138").Value = 7762.799999999999
$ws.Range("L138").Value = 15230.4228
$ws.Range("M138").Value = -2622.799999999999
$ws.Range("N138").Value = -25510.4228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 7000
$ws.Range("J14").Value = 7000
$ws.Range("L14").Value = 7000
$ws.Range("N14").Value = -7350
$ws.Range("H94").Value = 34886.668
$ws.Range("J94").Value = 34886.668
$ws.Range("L94").Value = 34886.668
$ws.Range("N94").Value = -36688.668
$ws.Range("H97").Value = 679.94116
$ws.Range("I97").Value = 630.6
$ws.Range("J97").Value = 1050
$ws.Range("K97").Value = 630.6
$ws.Range("L97").Value = 1050
$ws.Range("M97").Value = -134.6
$ws.Range("N97").Value = -2042
$ws.Range("H122").Value = 2049
$ws.Range("I122").Value = 2122.4
$ws.Range("K122").Value = 6367.200000000001
$ws.Range("M122").Value = -3917.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3718.8
$ws.Range("I134").Value = 3974.25
$ws.Range("K134").Value = 11922.75
$ws.Range("M134").Value = -9387.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 687.4
$ws.Range("I22").Value = 687
$ws.Range("J22").Value = 687.8
$ws.Range("K22").Value = 687
$ws.Range("L22").Value = 687.8
$ws.Range("M22").Value = -337
$ws.Range("N22").Value = -1387.8
$ws.Range("H54").Value = 39999
$ws.Range("J54").Value = 39999
$ws.Range("L54").Value = 39999
$ws.Range("N54").Value = -41315
$ws.Range("H94").Value = 118601
$ws.Range("J94").Value = 10202
$ws.Range("L94").Value = 10202
$ws.Range("N94").Value = -11104
$ws.Range("H132").Value = 1728.5
$ws.Range("I132").Value = 1640.45
$ws.Range("K132").Value = 4921.35
$ws.Range("M132").Value = -2391.35

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 276437.5
$ws.Range("J2").Value = 277800
$ws.Range("L2").Value = 1666800
$ws.Range("N2").Value = -1667026
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = 0
$ws.Range("H51").Value = 492.8889
$ws.Range("I51").Value = 429.625
$ws.Range("J51").Value = 999
$ws.Range("K51").Value = 1288.875
$ws.Range("L51").Value = 2997
$ws.Range("M51").Value = -828.875
$ws.Range("N51").Value = -3917
$ws.Range("H54").Value = 5148.3335
$ws.Range("J54").Value = 5148.3335
$ws.Range("L54").Value = 15445.0005
$ws.Range("N54").Value = -16563.0005
$ws.Range("H55").Value = 499.5
$ws.Range("I55").Value = 499.5
$ws.Range("K55").Value = 1498.5
$ws.Range("M55").Value = -1321.5
$ws.Range("H86").Value = 574.5714
$ws.Range("I86").Value = 713
$ws.Range("J86").Value = 390
$ws.Range("K86").Value = 2139
$ws.Range("L86").Value = 1170
$ws.Range("M86").Value = -953
$ws.Range("N86").Value = -3542
$ws.Range("H89").Value = 574.5714
$ws.Range("I89").Value = 713
$ws.Range("J89").Value = 390
$ws.Range("K89").Value = 6417
$ws.Range("L89").Value = 3510
$ws.Range("M89").Value = -489
$ws.Range("N89").Value = -15366
$ws.Range("H98").Value = 1216
$ws.Range("J98").Value = 1259.4
$ws.Range("L98").Value = 3778.2
$ws.Range("N98").Value = -6774.200000000001
$ws.Range("H107").Value = 562.6667
$ws.Range("J107").Value = 795
$ws.Range("L107").Value = 2385
$ws.Range("N107").Value = -6225

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3014.1765
$ws.Range("I102").Value = 2961.8
$ws.Range("K102").Value = 2961.8
$ws.Range("M102").Value = -1339.8
$ws.Range("H122").Value = 2488
$ws.Range("I122").Value = 2384.7144
$ws.Range("J122").Value = 2849.5
$ws.Range("K122").Value = 7154.1432
$ws.Range("L122").Value = 8548.5
$ws.Range("M122").Value = -4704.1432
$ws.Range("N122").Value = -13448.5
$ws.Range("H132").Value = 2381.2222
$ws.Range("I132").Value = 2381.2222
$ws.Range("K132").Value = 7143.6666
$ws.Range("M132").Value = -4613.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7520.6313
$ws.Range("I7").Value = 4300
$ws.Range("J7").Value = 8379.467000000001
$ws.Range("K7").Value = 4300
$ws.Range("L7").Value = 8379.467000000001
$ws.Range("M7").Value = -4188
$ws.Range("N7").Value = -8603.467000000001
$ws.Range("H16").Value = 849.6667
$ws.Range("I16").Value = 849.6667
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 849.6667
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = -679.6667
$ws.Range("H55").Value = 569.13336
$ws.Range("J55").Value = 498.33334
$ws.Range("L55").Value = 498.33334
$ws.Range("N55").Value = -844.33334
$ws.Range("H57").Value = 9999
$ws.Range("J57").Value = 9999
$ws.Range("L57").Value = 9999
$ws.Range("N57").Value = -11131
$ws.Range("H93").Value = 3075.8333
$ws.Range("I93").Value = 2633.3333
$ws.Range("K93").Value = 2633.3333
$ws.Range("M93").Value = -1385.3333
$ws.Range("H122").Value = 5149.6
$ws.Range("I122").Value = 4241.8945
$ws.Range("J122").Value = 6717.4546
$ws.Range("K122").Value = 12725.6835
$ws.Range("L122").Value = 20152.3638
$ws.Range("M122").Value = -10275.6835
$ws.Range("N122").Value = -25052.3638
$ws.Range("H126").Value = 7520.6313
$ws.Range("I126").Value = 4300
$ws.Range("J126").Value = 8379.467000000001
$ws.Range("K126").Value = 12900
$ws.Range("L126").Value = 25138.401
$ws.Range("M126").Value = -10430
$ws.Range("N126").Value = -30078.401
$ws.Range("H132").Value = 6767.1
$ws.Range("I132").Value = 6962.8887
$ws.Range("J132").Value = 5005
$ws.Range("K132").Value = 20888.6661
$ws.Range("L132").Value = 15015
$ws.Range("M132").Value = -18358.6661
$ws.Range("N132").Value = -20075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19993.25
$ws.Range("I41").Value = 19989
$ws.Range("J41").Value = 19994.666
$ws.Range("K41").Value = 19989
$ws.Range("L41").Value = 19994.666
$ws.Range("M41").Value = -19599
$ws.Range("N41").Value = -20774.666
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = ""
$ws.Range("N110").Value = 0
$ws.Range("H119").Value = 84999
$ws.Range("J119").Value = 84999
$ws.Range("L119").Value = 84999
$ws.Range("N119").Value = -94675
$ws.Range("H122").Value = 1940.091
$ws.Range("I122").Value = 1940.091
$ws.Range("K122").Value = 5820.272999999999
$ws.Range("M122").Value = -3370.272999999999
